# Apply the updated crypto price/volume snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.896.73'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.740.52'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('D5').Value = "'" + '248.22'
$ws.Range('E5').Value = '  +5.62%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = "'" + '0.5091'
$ws.Range('E7').Value = '  -3.71%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'" + '0.2732'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'" + '0.06179'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D10').Value = '1.740.75'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'" + '0.07232'
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = "'" + '15.10'
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = "'" + '0.6467'
$ws.Range('E13').Value = '  +0.31%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'" + '4.625'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'" + '77.60'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('B16').Value = 'Dai'
$ws.Range('C16').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D16').Value = "'" + '1.000'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = "'" + '1.000'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '25.922.70'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = "'" + '11.81'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'" + '0.000006798'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '1.963.99'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'" + '4.273'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = "'" + '8.639'
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = "'" + '5.387'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'" + '136.16'
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = "'" + '1.500'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'" + '15.22'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = "'" + '1.772'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'" + '105.30'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'" + '3.903'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = "'" + '0.08226'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'" + '3.632'
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'" + '0.04674'
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = "'" + '2.656'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'" + '0.9956'
$ws.Range('E35').Value = '  -0.89%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'" + '0.6242'
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = "'" + '2.729'
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'" + '0.01603'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'" + '1.921'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = "'" + '0.9998'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').Value = "'" + '99.87'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'" + '0.7600'
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = "'" + '0.3845'
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = "'" + '4.990'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = "'" + '0.1129'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = "'" + '6.274'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'" + '55.36'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'" + '0.05231'
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = "'" + '30.63'
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'" + '7.518'
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = "'" + '0.3410'
$ws.Range('E51').Value = '  -1.27%  '
